$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 31 (the lightgbm v2data_final submission): its "Name Ramp" note
# ("version_1_3") is cleared out and "Hand in" flips from TRUE to FALSE
# now that it has been superseded.
$ws.Cells.Item(31, 3).ClearContents()

$d31 = $ws.Cells.Item(31, 4)
$d31.Formula = "=""FALSE"""
$d31.Copy()
$d31.PasteSpecial(-4163) | Out-Null

# Row 32: new submission entry for the voting regressor.
$ws.Cells.Item(32, 1).Value = $ws.Cells.Item(31, 1).Value2
$ws.Cells.Item(32, 2).Value = "221127_voting_regressor"
$ws.Cells.Item(32, 3).Value = "last_but_not_best"

$d32 = $ws.Cells.Item(32, 4)
$d32.Formula = "=""TRUE"""
$d32.Copy()
$d32.PasteSpecial(-4163) | Out-Null

$ws.Cells.Item(32, 5).Value = "Maria"

$excel.CutCopyMode = 0

# Leave the selection where the user ended up after filling in the row.
$ws.Range("C34").Select() | Out-Null
